$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was inserted before the existing row 254,
# pushing every following record down by one row (old row 254 -> new
# row 255, ..., old row 323 -> new row 324). Insert a blank row at 254
# first so the rest of the table shifts down intact, then populate the
# new row with its own data.
$ws.Rows.Item(254).Insert()

$ws.Cells.Item(254, 1).Value = 8
$ws.Cells.Item(254, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(254, 3).Value = "Coquimbo"
$ws.Cells.Item(254, 4).Value = 44841
$ws.Cells.Item(254, 5).Value = 4
$ws.Cells.Item(254, 6).Value = 100112012
$ws.Cells.Item(254, 7).Value = "Espinaca"
$ws.Cells.Item(254, 8).Value = "Sin especificar"
$ws.Cells.Item(254, 9).Value = "Primera"
$ws.Cells.Item(254, 10).Value = 3000
$ws.Cells.Item(254, 11).Value = 450
$ws.Cells.Item(254, 12).Value = 500
$ws.Cells.Item(254, 13).Value = 475
$ws.Cells.Item(254, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(254, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(254, 16).Value = 950
$ws.Cells.Item(254, 17).Value = 0.5
$ws.Cells.Item(254, 18).Value = "Hortaliza"
